$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries & provincias update (data refreshed, table re-sorted by Casos totales desc) ---
# Rows 13/14 swap rank: Belgica now outranks Rusia
$ws.Range("A13").Value = "Belgica"   # was Rusia
$ws.Range("A14").Value = "Rusia"     # was Belgica

# Rows 40/41/42 rotate rank: Indonesia now outranks Filipinas and Singapur
$ws.Range("A40").Value = "Indonesia" # was Filipinas
$ws.Range("A41").Value = "Filipinas" # was Singapur
$ws.Range("A42").Value = "Singapur"  # was Indonesia

# --- Updated statistics ---
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Alemania (row 8)
$ws.Range("D8").Value = 85400
$ws.Range("E8").Value = 51645

# Belgica (row 13)
$ws.Range("B13").Value = 37183
$ws.Range("C13").Value = 1045
$ws.Range("D13").Value = 8348
$ws.Range("E13").Value = 23382
$ws.Range("F13").Value = 1119
$ws.Range("G13").Value = 290
$ws.Range("H13").Value = 5453

# Rusia (row 14)
$ws.Range("B14").Value = 36793
$ws.Range("C14").Value = 4785
$ws.Range("D14").Value = 3057
$ws.Range("E14").Value = 33423
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 40
$ws.Range("H14").Value = 313

# Austria (row 20)
$ws.Range("B20").Value = 14619
$ws.Range("C20").Value = 24
$ws.Range("E20").Value = 3974
$ws.Range("F20").Value = 208

# Rumania (row 31)
$ws.Range("E31").Value = 6142
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = 417

# Indonesia (row 40)
$ws.Range("B40").Value = 6248
$ws.Range("C40").Value = 325
$ws.Range("D40").Value = 631
$ws.Range("E40").Value = 5082
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 15
$ws.Range("H40").Value = 535

# Filipinas (row 41)
$ws.Range("B41").Value = 6087
$ws.Range("C41").Value = 209
$ws.Range("D41").Value = 516
$ws.Range("E41").Value = 5174
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 10
$ws.Range("H41").Value = 397

# Singapur (row 42)
$ws.Range("B42").Value = 5992
$ws.Range("C42").Value = 942
$ws.Range("D42").Value = 708
$ws.Range("E42").Value = 5273
$ws.Range("F42").Value = 22
$ws.Range("H42").Value = 11

# Marruecos (row 57)
$ws.Range("B57").Value = 2670
$ws.Range("C57").Value = 106
$ws.Range("D57").Value = 298
$ws.Range("E57").Value = 2235
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 137

# Banglades (row 61)
$ws.Range("B61").Value = 2144
$ws.Range("C61").Value = 306
$ws.Range("D61").Value = 66
$ws.Range("E61").Value = 1994
$ws.Range("G61").Value = 9
$ws.Range("H61").Value = 84

# Barein (row 65)
$ws.Range("B65").Value = 1744
$ws.Range("C65").Value = 4
$ws.Range("D65").Value = 726
$ws.Range("E65").Value = 1011

# Eslovaquia (row 79)
$ws.Range("F79").Value = 10

# Hong Kong (row 80)
$ws.Range("B80").Value = 1024
$ws.Range("C80").Value = 2
$ws.Range("D80").Value = 568
$ws.Range("E80").Value = 452

# Libano (row 92)
$ws.Range("B92").Value = 672
$ws.Range("C92").Value = 4
$ws.Range("E92").Value = 557
